$wb = $excel.ActiveWorkbook

# Add the new "Weibull" worksheet as the last sheet (after "DSS")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Weibull"

# Populate the Weibull MLE results table
$ws.Range("A1").Value = "Weibull"
$ws.Range("B1").Value = "a (MLE)"
$ws.Range("C1").Value = "b (MLE)"
$ws.Range("D1").Value = "c (MLE)"
$ws.Range("A2").Value = "SYS1"
$ws.Range("B2").Value = 172.52600000000001
$ws.Range("C2").Value = 0.00069605700000000003
$ws.Range("D2").Value = 0.67673899999999998
$ws.Range("A3").Value = "SYS2"
$ws.Range("B3").Value = 139.953
$ws.Range("C3").Value = 0.000072103799999999999
$ws.Range("D3").Value = 0.82241600000000004
$ws.Range("A4").Value = "SYS3"
$ws.Range("B4").Value = 281.51100000000002
$ws.Range("C4").Value = 0.00018371500000000001
$ws.Range("D4").Value = 0.91422199999999998
$ws.Range("A5").Value = "CSR1"
$ws.Range("B5").Value = -32.572499999999998
$ws.Range("C5").Value = -0.366035
$ws.Range("D5").Value = 0.168351
$ws.Range("A6").Value = "CSR2"
$ws.Range("B6").Value = 137.63
$ws.Range("C6").Value = 0.000171457
$ws.Range("D6").Value = 0.85021500000000005
$ws.Range("A7").Value = "CSR3"
$ws.Range("B7").Formula = "=-6.69106*10^16"
$ws.Range("C7").Formula = "=-3.67959*10^10"
$ws.Range("D7").Value = -6.06487
$ws.Range("A8").Value = "SS3"
$ws.Range("B8").Value = 518.29600000000005
$ws.Range("C8").Value = 0.0000488617
$ws.Range("D8").Value = 0.88542600000000005
$ws.Range("A9").Value = "J1"
$ws.Range("B9").Value = 374.29899999999998
$ws.Range("C9").Value = 0.00708093
$ws.Range("D9").Value = 0.999274
$ws.Range("A10").Value = "J2"
$ws.Range("B10").Value = 4140.51
$ws.Range("C10").Value = 0.00030728199999999998
$ws.Range("D10").Value = 1.0009300000000001
$ws.Range("A11").Value = "J3"
$ws.Range("B11").Value = 529.73099999999999
$ws.Range("C11").Value = 0.026499700000000001
$ws.Range("D11").Value = 1.0979399999999999
$ws.Range("A12").Value = "J4"
$ws.Range("B12").Value = 1170.8800000000001
$ws.Range("C12").Value = 0.0015353000000000001
$ws.Range("D12").Value = 0.99515200000000004
$ws.Range("A13").Value = "J5"
$ws.Range("B13").Value = -28478.799999999999
$ws.Range("C13").Value = -0.00017540399999999999
$ws.Range("D13").Value = 0.99993100000000001
$ws.Range("A14").Value = "S2"
$ws.Range("B14").Value = 67.207099999999997
$ws.Range("C14").Value = 0.000446465
$ws.Range("D14").Value = 0.70719299999999996
$ws.Range("A15").Value = "S2IF"
$ws.Range("B15").Value = 67.207099999999997
$ws.Range("C15").Value = 0.000446465
$ws.Range("D15").Value = 0.70719299999999996
$ws.Range("A16").Value = "S2FC"
$ws.Range("B16").Value = 54.558799999999998
$ws.Range("C16").Value = 0.208235
$ws.Range("D16").Value = 0.74129299999999998
$ws.Range("A17").Value = "S27"
$ws.Range("B17").Value = 45.578000000000003
$ws.Range("C17").Value = 0.00018128400000000001
$ws.Range("D17").Value = 0.84483600000000003
$ws.Range("A18").Value = "SS1"
$ws.Range("B18").Value = 571.05799999999999
$ws.Range("C18").Value = 0.020326899999999998
$ws.Range("D18").Value = 1.06165
$ws.Range("A19").Value = "SS3"
$ws.Range("B19").Value = 518.29600000000005
$ws.Range("C19").Value = 0.0000488617
$ws.Range("D19").Value = 0.88542600000000005
$ws.Range("A20").Value = "SS4"
$ws.Range("B20").Value = 199.369
$ws.Range("C20").Value = 3.3154400000000002
$ws.Range("D20").Value = 0.026833699999999999
$ws.Range("A21").Value = "CDS"
$ws.Range("B21").Value = 529.73099999999999
$ws.Range("C21").Value = 0.026499700000000001
$ws.Range("D21").Value = 1.0979399999999999
$ws.Range("A22").Value = "DATA1"
$ws.Range("B22").Value = -51.443399999999997
$ws.Range("C22").Value = -0.020090199999999999
$ws.Range("D22").Value = 0.99924000000000002
$ws.Range("A23").Value = "DATA2"
$ws.Range("B23").Value = 1268.73
$ws.Range("C23").Value = 0.0040494800000000003
$ws.Range("D23").Value = 1.0001800000000001
$ws.Range("A24").Value = "DATA3"
$ws.Range("B24").Value = -94.670400000000001
$ws.Range("C24").Value = -0.0198009
$ws.Range("D24").Value = 0.99855799999999995
$ws.Range("A25").Value = "DATA4"
$ws.Range("B25").Value = 686.34699999999998
$ws.Range("C25").Value = 0.0342043
$ws.Range("D25").Value = 1.01362
$ws.Range("A26").Value = "DATA5"
$ws.Range("B26").Value = 1425.15
$ws.Range("C26").Value = 0.021786900000000001
$ws.Range("D26").Value = 1.0055099999999999
$ws.Range("A27").Value = "DATA6"
$ws.Range("B27").Value = 5423.56
$ws.Range("C27").Value = 0.0688304
$ws.Range("D27").Value = 1.01288
$ws.Range("A28").Value = "DATA7"
$ws.Range("B28").Value = 821.12599999999998
$ws.Range("C28").Value = 0.0096719700000000002
$ws.Range("D28").Value = 1.1399699999999999
$ws.Range("A29").Value = "DATA8"
$ws.Range("B29").Value = 1022.44
$ws.Range("C29").Value = 0.0057271500000000003
$ws.Range("D29").Value = 1.4871399999999999
$ws.Range("A30").Value = "DATA9"
$ws.Range("B30").Value = -56.566600000000001
$ws.Range("C30").Value = -0.0034130800000000002
$ws.Range("D30").Value = 0.99691099999999999
$ws.Range("A31").Value = "DATA10"
$ws.Range("B31").Value = -659.71100000000001
$ws.Range("C31").Value = -0.016404100000000001
$ws.Range("D31").Value = 0.997166
$ws.Range("A32").Value = "DATA11"
$ws.Range("E32").Value = "diverging "
$ws.Range("A33").Value = "DATA12"
$ws.Range("E33").Value = "diverging "
$ws.Range("A34").Value = "DATA13"
$ws.Range("E34").Value = "diverging "
$ws.Range("A35").Value = "DATA14"
$ws.Range("B35").Value = 266
$ws.Range("C35").Value = 11.184699999999999
$ws.Range("D35").Value = -1.23759

# Match the original selection state (B1 selected)
$ws.Range("B1").Select()
